# Generate Report for Handoff
#
# Updates the localization-status report to reflect a newer handoff run for
# the four "Ready for handoff" files (61fb34a1…, 8dda6b2d…, b4ee3680…,
# b8172df9…): their handoff timestamps advance and their Priority flips
# from "low" to "ht".

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G) for rows 4-7
# advances from 06:33:20 to 06:33:37.
$overview.Range("G4:G7").Value = "2016-09-03 06:33:37"

# zh-cn sheet: Priority (column E) goes from "low" to "ht" and the
# "Latest Handoff Datetime" (column H) advances from 06:33:16 to 06:33:32.
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4:H7").Value = "2016-09-03 06:33:32"

# de-de sheet: Priority (column E) also goes from "low" to "ht", and its
# "Latest Handoff Datetime" (column H) mirrors the Overview sheet's
# "Latest HO Xliff Generate Date", advancing from 06:33:20 to 06:33:37.
$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4:H7").Value = "2016-09-03 06:33:37"
